# "3. Kapitel wurde geschrieben"
#
# The two existing paragraphs get their single run split into several
# runs, with w:proofErr spellStart/spellEnd markers bracketing each German
# word (as Word's background spell-checker would leave behind after an
# edit), and a brand-new third paragraph ("Ozlem schrieb das dritte
# Kapitel!") is appended in the same style, preceded by a blank paragraph,
# mirroring the Mehmet/Akkaya pattern already in the document.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Split-SentenceXml($name, $verb, $ordinal) {
    $p  = '<w:p>'
    $p += '<w:r><w:t xml:space="preserve">' + $name + ' </w:t></w:r>'
    $p += '<w:proofErr w:type="spellStart"/><w:r><w:t>' + $verb + '</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    $p += '<w:r><w:t xml:space="preserve"> das </w:t></w:r>'
    $p += '<w:proofErr w:type="spellStart"/><w:r><w:t>' + $ordinal + '</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    $p += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
    $p += '<w:proofErr w:type="spellStart"/><w:r><w:t>Kapitel</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    $p += '<w:r><w:t>!</w:t></w:r>'
    $p += '</w:p>'
    return $p
}

# 1) Re-split "Mehmet schrieb das erste Kapitel!" (paragraph 1) into runs
#    with proofErr wraps around the German words.
$mehmetXml = $pkgHeader + (Split-SentenceXml "Mehmet" "schrieb" "erste") + $pkgFooter
$null = $d.Paragraphs.Item(1).Range.InsertXML($mehmetXml)

# 2) Re-split "Akkaya schrieb das zweite Kapitel!" (paragraph 3) likewise.
$akkayaXml = $pkgHeader + (Split-SentenceXml "Akkaya" "schrieb" "zweite") + $pkgFooter
$null = $d.Paragraphs.Item(3).Range.InsertXML($akkayaXml)

# 3) Append a blank paragraph and the new "Ozlem schrieb das dritte
#    Kapitel!" paragraph (same run-split pattern) at the end of the
#    document, ahead of the final, already-existing blank paragraph.
#    Replacing the full range of that final (empty) paragraph with three
#    paragraphs (blank, Ozlem sentence, blank) expands it in place: the
#    last of the three inherits the original trailing paragraph's
#    identity, so the document keeps exactly one blank paragraph at the
#    very end instead of growing an extra one.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$ozlemXml = $pkgHeader + '<w:p/>' + (Split-SentenceXml "Ozlem" "schrieb" "dritte") + '<w:p/>' + $pkgFooter
$null = $d.Range($lastPara.Range.Start, $lastPara.Range.End).InsertXML($ozlemXml)
